$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 0.105
$ws.Range("E2").Value = 0.0247
$ws.Range("G2").Value = 0.0619669277632724
$ws.Range("H2").Value = 0.0619669277632724
$ws.Range("I2").Value = 0.06544821583986075
$ws.Range("J2").Value = 0.05208203091481877
$ws.Range("K2").Value = 6.1
$ws.Range("L2").Value = 0.05308964316797214
$ws.Range("M2").Value = 3.84
$ws.Range("N2").Value = 0.05189189189189189
$ws.Range("O2").Value = 0.6295081967213115
$ws.Range("P2").Value = 3.84
$ws.Range("Q2").Value = 0.05189189189189189
$ws.Range("R2").Value = 0.6295081967213115
$ws.Range("U2").Value = 11.62
$ws.Range("V2").Value = 0.157027027027027
$ws.Range("W2").Value = 0.2579789124887326
$ws.Range("X2").Value = 0.04672939245001846
$ws.Range("Y2").Value = 0.2112495200387141
$ws.Range("Z2").Value = 5.302261190586064
$ws.Range("AA2").Value = 0.3168470612550844
$ws.Range("AB2").Value = 0.04461634212803336
$ws.Range("AC2").Value = 0.272230719127051
$ws.Range("AD2").Value = 5.88
$ws.Range("AE2").Value = 0
$ws.Range("AF2").Value = 5.88
$ws.Range("AG2").Value = -5.740000000000001
$ws.Range("AH2").Value = 0.07361041562343515
$ws.Range("AI2").Value = 0.1611842105263158
$ws.Range("AJ2").Value = -0.08409024318781132
$ws.Range("AK2").Value = -0.2308930008045053
$ws.Range("AL2").Value = 0.08699999999999999
$ws.Range("AM2").Value = -0.109
$ws.Range("AN2").Value = 0.7144592952612393
$ws.Range("AO2").Value = 86.43678160919541
$ws.Range("AP2").Value = -0.6974483596597814
$ws.Range("AQ2").Value = -68.99082568807339

# Row 3
$ws.Range("G3").Value = 0.04951100244498777
$ws.Range("H3").Value = 0.04951100244498777
$ws.Range("I3").Value = 0.05073349633251834
$ws.Range("J3").Value = 0.04058679706601467
$ws.Range("K3").Value = 3.29
$ws.Range("L3").Value = 0.0402200488997555
$ws.Range("M3").Value = 1.25
$ws.Range("N3").Value = 0.03720238095238095
$ws.Range("O3").Value = 0.3799392097264437
$ws.Range("P3").Value = 1.25
$ws.Range("Q3").Value = 0.03720238095238095
$ws.Range("R3").Value = 0.3799392097264437
$ws.Range("U3").Value = 9.210000000000001
$ws.Range("V3").Value = 0.2741071428571429
$ws.Range("W3").Value = 0.3466807165437302
$ws.Range("X3").Value = 0.0481976031932337
$ws.Range("Y3").Value = 0.2984831133504965
$ws.Range("Z3").Value = 10.99462365591398
$ws.Range("AA3").Value = 0.4462365591397849
$ws.Range("AB3").Value = 0.04476449848052817
$ws.Range("AC3").Value = 0.4014720606592567
$ws.Range("AD3").Value = 4.59
$ws.Range("AE3").Value = 0
$ws.Range("AF3").Value = 4.59
$ws.Range("AG3").Value = -4.620000000000001
$ws.Range("AH3").Value = 0.1201885310290652
$ws.Range("AI3").Value = 0.2654713707345286
$ws.Range("AJ3").Value = -0.1594202898550725
$ws.Range("AK3").Value = -0.5717821782178221
$ws.Range("AL3").Value = 0.055
$ws.Range("AM3").Value = 0.047
$ws.Range("AN3").Value = 1
$ws.Range("AO3").Value = 75.45454545454547
$ws.Range("AP3").Value = -1.006535947712418
$ws.Range("AQ3").Value = 88.29787234042554

# Row 4
$ws.Range("B4").Value = "Taiming Assurance Broker Co.,Ltd. (GTSM:5878)"
$ws.Range("D4").Value = 0.105
$ws.Range("E4").Value = 0.0247
$ws.Range("G4").Value = 0.09274924471299092
$ws.Range("H4").Value = 0.09274924471299092
$ws.Range("I4").Value = 0.1018126888217523
$ws.Range("J4").Value = 0.08058976213778138
$ws.Range("K4").Value = 2.81
$ws.Range("L4").Value = 0.08489425981873111
$ws.Range("M4").Value = 2.59
$ws.Range("N4").Value = 0.06410891089108911
$ws.Range("O4").Value = 0.9217081850533807
$ws.Range("P4").Value = 2.59
$ws.Range("Q4").Value = 0.06410891089108911
$ws.Range("R4").Value = 0.9217081850533807
$ws.Range("U4").Value = 2.41
$ws.Range("V4").Value = 0.05965346534653466
$ws.Range("W4").Value = 0.1692771084337349
$ws.Range("X4").Value = 0.04526118170680322
$ws.Range("Y4").Value = 0.1240159267269317
$ws.Range("Z4").Value = 2.326071679550246
$ws.Range("AA4").Value = 0.1874575633703839
$ws.Range("AB4").Value = 0.04446818577553856
$ws.Range("AC4").Value = 0.1429893775948454
$ws.Range("AD4").Value = 1.29
$ws.Range("AE4").Value = 0
$ws.Range("AF4").Value = 1.29
$ws.Range("AG4").Value = -1.12
$ws.Range("AH4").Value = 0.03094267210362197
$ws.Range("AI4").Value = 0.0672225117248567
$ws.Range("AJ4").Value = -0.0285132382892057
$ws.Range("AK4").Value = -0.06674612634088202
$ws.Range("AL4").Value = 0.032
$ws.Range("AM4").Value = -0.156
$ws.Range("AN4").Value = 0.3543956043956044
$ws.Range("AO4").Value = 105.3125
$ws.Range("AP4").Value = -0.3076923076923077
$ws.Range("AQ4").Value = -21.6025641025641

